$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.857.27"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.357.59"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.97"
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("E6").Value = "  -1.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.35"
$ws.Range("E7").Value = "  -1.55%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  +1.78%  "

$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.69"
$ws.Range("E11").Value = "  +6.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "34.05"
$ws.Range("E12").Value = "  +5.65%  "

$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.17"
$ws.Range("E15").Value = "  -2.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.909"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.358.08"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.900.10"
$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.73"
$ws.Range("E20").Value = "  +0.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.52"
$ws.Range("E21").Value = "  -3.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "252.74"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("E24").Value = "  +2.88%  "

$ws.Range("E25").Value = "  -5.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.49"
$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("E27").Value = "  -2.47%  "

$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "176.28"
$ws.Range("E29").Value = "  +0.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.26"
$ws.Range("E30").Value = "  -2.28%  "

$ws.Range("E31").Value = "  +0.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.133"
$ws.Range("E32").Value = "  -2.40%  "

$ws.Range("E33").Value = "  -1.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.05"
$ws.Range("E34").Value = "  -3.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.32"
$ws.Range("E35").Value = "  -1.94%  "

$ws.Range("E36").Value = "  +1.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.59"
$ws.Range("E37").Value = "  +4.46%  "

$ws.Range("E38").Value = "  +1.47%  "

$ws.Range("E39").Value = "  -1.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.51"
$ws.Range("E40").Value = "  +16.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.72"
$ws.Range("E41").Value = "  +10.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.95"
$ws.Range("E42").Value = "  +3.97%  "

$ws.Range("E43").Value = "  -6.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.02"
$ws.Range("E44").Value = "  -1.61%  "

$ws.Range("E45").Value = "  -2.49%  "

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("E47").Value = "  -0.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.45"
$ws.Range("E48").Value = "  -1.96%  "

$ws.Range("E49").Value = "  -1.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "97.97"
$ws.Range("E50").Value = "  -2.40%  "

$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000212"
$ws.Range("E51").Value = "  +14.15%  "
